# KIBON-917: Lastenausgleich Excel Export: Mit Totals
#
# Adds a set of workbook-level defined names that point at the single
# "input row" (row 8) of the Data sheet, and turns the previously-empty
# totals row (row 9) into SUM() formulas over those named ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Workbook-level defined names, one per aggregated column in row 8.
$wb.Names.Add("totalBelegung", "=Data!`$D`$8")
$wb.Names.Add("totalGutscheine", "=Data!`$E`$8")
$wb.Names.Add("kostenProHundertProzentPlatz", "=Data!`$F`$8")
$wb.Names.Add("selbstbehaltGemeinde", "=Data!`$G`$8")
$wb.Names.Add("eingabeLastenausgleich", "=Data!`$H`$8")

# Totals row: sum each named range into the corresponding column of row 9.
$ws.Range("D9").Formula = "=SUM(totalBelegung)"
$ws.Range("E9").Formula = "=SUM(totalGutscheine)"
$ws.Range("F9").Formula = "=SUM(kostenProHundertProzentPlatz)"
$ws.Range("G9").Formula = "=SUM(selbstbehaltGemeinde)"
$ws.Range("H9").Formula = "=SUM(eingabeLastenausgleich)"
